$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 2019
$ws.Range("A5").Value = 2018
$ws.Range("A6").Value = 2017
$ws.Range("A7").Value = 2016
$ws.Range("A8").Value = 2015
$ws.Range("A9").Value = 2014

$ws.Range("A10").Select()
